$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update text content (column B) ---
$ws.Range("B6").Value = "Thank you, Cuong! I appreciate your willingness to help. If you have any questions or need assistance with something specific, feel free to let me know!"
$ws.Range("B10").Value = "Hello Cuong! I'm here to assist you as well. How can I help you today?"
$ws.Range("B11").Value = "Hello! I appreciate your offer to help. I'm here to assist you with any questions or information you need. What would you like to know or discuss today?"
$ws.Range("B12").Value = "Thank you, Cuong! I appreciate your willingness to assist. I’m here to provide information and answer any questions you might have. Is there a specific topic or question you’d like to discuss?"

# --- Update Response_Time values (column C) ---
$ws.Range("C2").Value = 1.275030851364136
$ws.Range("C3").Value = 1.213689804077148
$ws.Range("C4").Value = 1.02789568901062
$ws.Range("C5").Value = 1.816789388656616
$ws.Range("C8").Value = 0.7870028018951416
$ws.Range("C9").Value = 0.8151404857635498
$ws.Range("C10").Value = 1.332005023956299
$ws.Range("C11").Value = 1.669925451278687
